# Gliederung.docx update
#  1. "Skalierung" -> "Skalierbarkeit" (both the outline entry and the
#     matching heading further down in the "Auswertung" section).
#  2. In the "Auswertung" (Experimente) section, the "Testplan Skalierung"
#     block used to start with a stray "Datendurchsatz" bullet copied from
#     the previous block. That bullet is removed from there and, instead,
#     two bullets - "Median, Mittelwert, Abweichung, Minimum, Maximum und
#     Schwankung der Latenz" and "Datendurchsatz" - are appended under
#     "Testplan Last", mirroring the bullets already used under
#     "Testplan Performanz".

$d = $word.ActiveDocument

# --- 1) Global rename: Skalierung -> Skalierbarkeit -------------------
$d.Content.Find.Execute("Skalierung", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Skalierbarkeit", 2) | Out-Null

# --- 2) Work only inside the "Experimente" section so the outline copies
#        of these headings (earlier in the document) are left untouched.
$sectionStart = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "Experimente") {
        $sectionStart = $i
        break
    }
}

$testplanLastIndex = -1
$oldDatendurchsatzIndex = -1
for ($i = $sectionStart; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "Testplan Last" -and $testplanLastIndex -eq -1) {
        $testplanLastIndex = $i
    }
    if ($txt -eq "Datendurchsatz" -and $oldDatendurchsatzIndex -eq -1) {
        $oldDatendurchsatzIndex = $i
    }
}

# --- 3) Append the two new bullets right after "Testplan Last" ---------
$testplanLast = $d.Paragraphs.Item($testplanLastIndex)
$testplanLast.Range.InsertParagraphAfter()

$medianPara = $d.Paragraphs.Item($testplanLastIndex + 1)
$medianPara.Range.Text = "Median, Mittelwert, Abweichung, Minimum, Maximum und Schwankung der Latenz"
$medianPara.Range.ListFormat.ListLevelNumber = 4

$medianPara.Range.InsertParagraphAfter()
$datenPara = $d.Paragraphs.Item($testplanLastIndex + 2)
$datenPara.Range.Text = "Datendurchsatz"
$datenPara.Range.ListFormat.ListLevelNumber = 4

# --- 4) Remove the old, now-duplicated "Datendurchsatz" bullet that used
#        to sit right before "Testplan Skalierbarkeit" in this section.
$d.Paragraphs.Item($oldDatendurchsatzIndex).Range.Delete()

Write-Output "Done."
